$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> Alvearie Team
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting all following rows up by one.
$ws1.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root element row: Short "Extension" -> "Language Rank"
$ws2.Range("K2").Value = "Language Rank"

# Root element row: Definition "An Extension" -> "Specify preferred order of language use (1 = highest)"
$ws2.Range("L2").Value = "Specify preferred order of language use (1 = highest)"
